$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "28.223.74", "  +3.43%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.814.60", "  +4.56%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.002", "  -0.43%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "328.89", "  +1.76%  "),
    @(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.9992", "  -0.50%  "),
    @(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4433", "  +4.20%  "),
    @(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.3708", "  +3.25%  "),
    @(9, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "44.89", "  +0.08%  "),
    @(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.07703", "  +4.57%  "),
    @(11, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.124", "  +0.71%  "),
    @(12, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.001", "  -0.37%  "),
    @(13, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "22.04", "  +2.72%  "),
    @(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "6.257", "  +3.20%  "),
    @(15, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.561", "  +5.51%  "),
    @(16, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.789.36", "  +2.97%  "),
    @(17, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "92.92", "  +10.25%  "),
    @(18, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.00001082", "  +2.55%  "),
    @(19, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.06530", "  +9.49%  "),
    @(20, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.9991", "  -0.52%  "),
    @(21, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "17.52", "  +4.57%  "),
    @(22, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.218", "  +3.47%  "),
    @(23, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "28.296.99", "  +3.58%  "),
    @(24, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "11.69", "  +3.80%  "),
    @(25, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "20.66", "  +4.03%  "),
    @(26, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.932", "  -19.56%  "),
    @(27, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "155.77", "  +4.69%  "),
    @(28, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.019.56", "  +4.25%  "),
    @(29, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "2.318", "  -0.48%  "),
    @(30, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "127.96", "  +1.98%  "),
    @(31, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "1.200", "  -4.65%  "),
    @(32, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "5.875", "  +6.13%  "),
    @(33, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.09215", "  +2.29%  "),
    @(34, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "3.672", "  -0.96%  "),
    @(35, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "13.13", "  +7.13%  "),
    @(36, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.02349", "  +3.48%  "),
    @(37, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.2167", "  +0.50%  "),
    @(38, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "5.173", "  +3.83%  "),
    @(39, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.6566", "  +2.63%  "),
    @(40, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.06189", "  +1.88%  "),
    @(41, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "1.198", "  +1.70%  "),
    @(42, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "8.084", "  +3.59%  "),
    @(43, "Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "0.9985", "  -0.55%  "),
    @(44, "WEMIXTOKEN", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "1.386", "  -1.73%  "),
    @(45, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "13.84", "  +3.37%  "),
    @(46, "Decentraland", "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana", "0.6075", "  +4.00%  "),
    @(47, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "3.757", "  +0.46%  "),
    @(48, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "2.036", "  +5.60%  "),
    @(49, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "126.54", "  +1.70%  "),
    @(50, "EOS", "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos", "1.153", "  +5.72%  "),
    @(51, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.06979", "  +2.65%  ")
)

$ws.Range("D2:D51").NumberFormat = "@"

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
